$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

# Column A: a date-like string ("08/06/2025") that must stay a literal text
# value (matching the existing A10:A37 cells) rather than being
# auto-converted by Excel into a date serial number. Temporarily force a
# text number format while the value is entered, then restore the cell to
# the default "Normal" style so it matches the unstyled cells above it.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "08/06/2025"
$cellA.Style = "Normal"

# Columns B-D: plain numeric values.
$ws.Cells.Item($row, 2).Value = 0.0004301200000000026
$ws.Cells.Item($row, 3).Value = 116246.6288477627
$ws.Cells.Item($row, 4).Value = 50
